$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 13.2188964569129
$ws.Range("D2").Value = 5.015407118611042
$ws.Range("E2").Value = 15.73769832470673
$ws.Range("F2").Value = 24.62405010847608
$ws.Range("G2").Value = 3.63384671211641
$ws.Range("I2").Value = 26.04689624712438
$ws.Range("K2").Value = 10.06641956127041
$ws.Range("L2").Value = 8.93886298257522
$ws.Range("O2").Value = 22.05706760004436

$ws.Range("B3").Value = 13.06304640938694
$ws.Range("D3").Value = 4.96388766446606
$ws.Range("E3").Value = 15.79668642953868
$ws.Range("F3").Value = 24.6481426277848
$ws.Range("G3").Value = 3.635699780579466
$ws.Range("I3").Value = 26.17763357169516
$ws.Range("K3").Value = 9.705362498731322
$ws.Range("L3").Value = 8.892993726301984
$ws.Range("O3").Value = 22.12565081014536

$ws.Range("B4").Value = 12.96856150794201
$ws.Range("D4").Value = 4.931554511267962
$ws.Range("E4").Value = 15.83508814700937
$ws.Range("F4").Value = 24.67002506591991
$ws.Range("G4").Value = 3.636898409538476
$ws.Range("I4").Value = 26.26262771819591
$ws.Range("K4").Value = 9.474879663967492
$ws.Range("L4").Value = 8.866270394642449
$ws.Range("O4").Value = 22.17302059992018

$ws.Range("B5").Value = 12.93040303106728
$ws.Range("D5").Value = 4.918208764351297
$ws.Range("E5").Value = 15.85128713316782
$ws.Range("F5").Value = 24.6807219118791
$ws.Range("G5").Value = 3.637402206188348
$ws.Range("I5").Value = 26.29845187796449
$ws.Range("K5").Value = 9.378833264946856
$ws.Range("L5").Value = 8.855751175386551
$ws.Range("O5").Value = 22.19364323933021

$ws.Range("B6").Value = 12.92408886911441
$ws.Range("D6").Value = 4.915982648273244
$ws.Range("E6").Value = 15.85401021196162
$ws.Range("F6").Value = 24.68260548565806
$ws.Range("G6").Value = 3.637486789474979
$ws.Range("I6").Value = 26.30447226356279
$ws.Range("K6").Value = 9.362759401268454
$ws.Range("L6").Value = 8.854027098772109
$ws.Range("O6").Value = 22.19714717657296

$ws.Range("B7").Value = 12.96804543889242
$ws.Range("D7").Value = 4.931375204682648
$ws.Range("E7").Value = 15.83530438390315
$ws.Range("F7").Value = 24.67016212680118
$ws.Range("G7").Value = 3.636905141719527
$ws.Range("I7").Value = 26.2631060420606
$ws.Range("K7").Value = 9.473592819672747
$ws.Range("L7").Value = 8.866127016720259
$ws.Range("O7").Value = 22.17329338792223

$ws.Range("B8").Value = 13.16493358065271
$ws.Range("D8").Value = 4.997792236410533
$ws.Range("E8").Value = 15.75758504658113
$ws.Range("F8").Value = 24.63088417689929
$ws.Range("G8").Value = 3.634473051982713
$ws.Range("I8").Value = 26.09099540043034
$ws.Range("K8").Value = 9.943801671320449
$ws.Range("L8").Value = 8.922753522388325
$ws.Range("O8").Value = 22.07962147177644

$ws.Range("B9").Value = 13.55878443625989
$ws.Range("D9").Value = 5.122212369834995
$ws.Range("E9").Value = 15.62244725683787
$ws.Range("F9").Value = 24.61021041783568
$ws.Range("G9").Value = 3.630184278959718
$ws.Range("I9").Value = 25.79088809081255
$ws.Range("K9").Value = 10.79274914496407
$ws.Range("L9").Value = 9.044832958264751
$ws.Range("O9").Value = 21.93780992520941

$ws.Range("B10").Value = 13.85041873830274
$ws.Range("D10").Value = 5.209724270509019
$ws.Range("E10").Value = 15.53362013391976
$ws.Range("F10").Value = 24.62944444661906
$ws.Range("G10").Value = 3.627323233819872
$ws.Range("I10").Value = 25.59311990932721
$ws.Range("K10").Value = 11.36802716303395
$ws.Range("L10").Value = 9.140693897241219
$ws.Range("O10").Value = 21.85934260727261

$ws.Range("B11").Value = 13.98306777537621
$ws.Range("D11").Value = 5.248617787201152
$ws.Range("E11").Value = 15.49546677039608
$ws.Range("F11").Value = 24.64566161926043
$ws.Range("G11").Value = 3.626083976857783
$ws.Range("I11").Value = 25.50806752761157
$ws.Range("K11").Value = 11.61856637665561
$ws.Range("L11").Value = 9.185516287736236
$ws.Range("O11").Value = 21.82926762818672

$ws.Range("B12").Value = 14.03325369623563
$ws.Range("D12").Value = 5.263208209515216
$ws.Range("E12").Value = 15.48134220877749
$ws.Range("F12").Value = 24.65287388075508
$ws.Range("G12").Value = 3.625623604606476
$ws.Range("I12").Value = 25.47656593949564
$ws.Range("K12").Value = 11.71178974779542
$ws.Range("L12").Value = 9.202652856602199
$ws.Range("O12").Value = 21.81868981877902

$ws.Range("B13").Value = 14.02244795304758
$ws.Range("D13").Value = 5.260072118815247
$ws.Range("E13").Value = 15.4843698218211
$ws.Range("F13").Value = 24.65127300516686
$ws.Range("G13").Value = 3.625722358535927
$ws.Range("I13").Value = 25.48331898430067
$ws.Range("K13").Value = 11.69178649238175
$ws.Range("L13").Value = 9.198955110867802
$ws.Range("O13").Value = 21.82093183174741

$ws.Range("B14").Value = 13.98719775374355
$ws.Range("D14").Value = 5.249820943481227
$ws.Range("E14").Value = 15.49429826017877
$ws.Range("F14").Value = 24.64623352781501
$ws.Range("G14").Value = 3.626045923484698
$ws.Range("I14").Value = 25.50546172904746
$ws.Range("K14").Value = 11.62626923629116
$ws.Range("L14").Value = 9.186922916838755
$ws.Range("O14").Value = 21.82838111828709

$ws.Range("B15").Value = 13.96559881839263
$ws.Range("D15").Value = 5.243523694678175
$ws.Range("E15").Value = 15.50042179036019
$ws.Range("F15").Value = 24.64328610193595
$ws.Range("G15").Value = 3.626245275086851
$ws.Range("I15").Value = 25.5191167112853
$ws.Range("K15").Value = 11.58592184128415
$ws.Range("L15").Value = 9.179573775221876
$ws.Range("O15").Value = 21.83304970608122

$ws.Range("B16").Value = 13.84174577894785
$ws.Range("D16").Value = 5.207163538684599
$ws.Range("E16").Value = 15.53615883484211
$ws.Range("F16").Value = 24.62853468179456
$ws.Range("G16").Value = 3.627405471010644
$ws.Range("I16").Value = 25.59877709636329
$ws.Range("K16").Value = 11.35142524125453
$ws.Range("L16").Value = 9.137788126128923
$ws.Range("O16").Value = 21.86142141198656

$ws.Range("B17").Value = 13.76573078155406
$ws.Range("D17").Value = 5.184618819302631
$ws.Range("E17").Value = 15.55865914287821
$ws.Range("F17").Value = 24.62139625683766
$ws.Range("G17").Value = 3.628133125433548
$ws.Range("I17").Value = 25.64890412162858
$ws.Range("K17").Value = 11.20467666986714
$ws.Range("L17").Value = 9.112457104458807
$ws.Range("O17").Value = 21.88026804324987

$ws.Range("B18").Value = 13.72200989780533
$ws.Range("D18").Value = 5.171565812819711
$ws.Range("E18").Value = 15.57181298647947
$ws.Range("F18").Value = 24.61799353369686
$ws.Range("G18").Value = 3.628557514558311
$ws.Range("I18").Value = 25.67819836290645
$ws.Range("K18").Value = 11.11922287254973
$ws.Range("L18").Value = 9.098002409968824
$ws.Range("O18").Value = 21.89163697775717

$ws.Range("B19").Value = 13.70720824162562
$ws.Range("D19").Value = 5.16713172419196
$ws.Range("E19").Value = 15.5763031360336
$ws.Range("F19").Value = 24.61696225168072
$ws.Range("G19").Value = 3.628702213482746
$ws.Range("I19").Value = 25.68819634513356
$ws.Range("K19").Value = 11.09011121355243
$ws.Range("L19").Value = 9.093128399218052
$ws.Range("O19").Value = 21.89557704614088

$ws.Range("B20").Value = 13.7738229254282
$ws.Range("D20").Value = 5.187027675309499
$ws.Range("E20").Value = 15.55624198506066
$ws.Range("F20").Value = 24.62208340059485
$ws.Range("G20").Value = 3.62805505904661
$ws.Range("I20").Value = 25.64352015173869
$ws.Range("K20").Value = 11.22040711274119
$ws.Range("L20").Value = 9.115141810433034
$ws.Range("O20").Value = 21.87820703232873

$ws.Range("B21").Value = 13.99755316139604
$ws.Range("D21").Value = 5.252835746244991
$ws.Range("E21").Value = 15.49137327037619
$ws.Range("F21").Value = 24.64768470009643
$ws.Range("G21").Value = 3.62595064321565
$ws.Range("I21").Value = 25.4989387148873
$ws.Range("K21").Value = 11.64555836822946
$ws.Range("L21").Value = 9.190452723870028
$ws.Range("O21").Value = 21.8261710531104

$ws.Range("B22").Value = 14.14349000522525
$ws.Range("D22").Value = 5.295039815315941
$ws.Range("E22").Value = 15.45086177984782
$ws.Range("F22").Value = 24.67065817915391
$ws.Range("G22").Value = 3.624627185384321
$ws.Range("I22").Value = 25.40856068210889
$ws.Range("K22").Value = 11.91378259779231
$ws.Range("L22").Value = 9.240619161809052
$ws.Range("O22").Value = 21.79689036060818

$ws.Range("B23").Value = 14.06564067549667
$ws.Range("D23").Value = 5.272590321423468
$ws.Range("E23").Value = 15.47231143892061
$ws.Range("F23").Value = 24.65782687442873
$ws.Range("G23").Value = 3.625328805219361
$ws.Range("I23").Value = 25.45642087337366
$ws.Range("K23").Value = 11.7715215956371
$ws.Range("L23").Value = 9.213761617087558
$ws.Range("O23").Value = 21.81208457794248

$ws.Range("B24").Value = 13.77016452260965
$ws.Range("D24").Value = 5.185938916821542
$ws.Range("E24").Value = 15.55733410264095
$ws.Range("F24").Value = 24.62177055835434
$ws.Range("G24").Value = 3.628090334008712
$ws.Range("I24").Value = 25.64595276316543
$ws.Range("K24").Value = 11.21329875893529
$ws.Range("L24").Value = 9.113927716201941
$ws.Range("O24").Value = 21.87913715294741

$ws.Range("B25").Value = 13.45166128002592
$ws.Range("D25").Value = 5.089213489086176
$ws.Range("E25").Value = 15.65716400986894
$ws.Range("F25").Value = 24.60975895486705
$ws.Range("G25").Value = 3.631293375162708
$ws.Range("I25").Value = 25.86807951736849
$ws.Range("K25").Value = 10.57134738271456
$ws.Range("L25").Value = 9.010681860186091
$ws.Range("O25").Value = 21.97166970345615
